# Update "想去人数" (people interested) counts in column F across the
# "展览" (sheet1), "演出" (sheet2), "本地生活" (sheet3) and "全部类型" (sheet4)
# worksheets to match the refreshed data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

function Set-FValues {
    param(
        [string]$SheetName,
        [hashtable]$Updates
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $Updates.Keys) {
        $ws.Range("F$row").Value = $Updates[$row]
    }
}

# 展览 (Exhibitions)
Set-FValues "展览" @{
    6  = 282
    7  = 13245
    8  = 77
    10 = 309
    11 = 4953
    13 = 3601
    14 = 48
    17 = 181
    20 = 48
    23 = 86
    24 = 112
    25 = 4486
    27 = 1960
    28 = 112
    29 = 276
    30 = 7121
    33 = 2143
    34 = 2065
    36 = 123
    37 = 1109
    39 = 6
    40 = 235
    42 = 1158
    43 = 9
    44 = 153
    45 = 1251
    46 = 1876
    47 = 82
}

# 演出 (Performances)
Set-FValues "演出" @{
    4 = 34
    8 = 134
}

# 本地生活 (Local life)
Set-FValues "本地生活" @{
    2 = 491
    3 = 664
    4 = 45
}

# 全部类型 (All types)
Set-FValues "全部类型" @{
    5  = 491
    6  = 664
    7  = 282
    8  = 13245
    10 = 309
    11 = 4953
    12 = 3601
    13 = 48
    15 = 181
    17 = 48
    21 = 86
    23 = 112
    24 = 4486
    26 = 1960
    27 = 112
    28 = 276
    29 = 7121
    33 = 2143
    34 = 2065
    36 = 123
    37 = 1109
    38 = 6
    39 = 235
    41 = 1158
    42 = 153
    44 = 1251
    45 = 1876
    46 = 82
}
